$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "Ba Form" conjugations
# Column width (matches diff col min=7 max=7 width~24.7109375)
$ws.Columns.Item(7).ColumnWidth = 23.996651785714285

# Copy existing cell formats into column G first, so the new cells reuse
# the same shared style records (fonts/alignment) as their row siblings
# instead of minting fresh ones.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("G2:G8").PasteSpecial(-4122)
$ws.Range("F20").Copy()
$ws.Range("G9:G88").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the Ba-Form values, in the same cell order the workbook was
# originally edited so new shared-string entries line up with the source.
$ws.Range("G2").Value = "言えば"
$ws.Range("G3").Value = "飲めば"
$ws.Range("G4").Value = "書けば"
$ws.Range("G5").Value = "話せば"
$ws.Range("G6").Value = "買えば"
$ws.Range("G1").Value = "Ba Form"
$ws.Range("G7").Value = "忘れれば"
$ws.Range("G8").Value = "読めば"
$ws.Range("G9").Value = "na"
$ws.Range("G10").Value = "na"
$ws.Range("G11").Value = "na"
$ws.Range("G12").Value = "na"
$ws.Range("G13").Value = "na"
$ws.Range("G14").Value = "na"
$ws.Range("G15").Value = "na"
$ws.Range("G16").Value = "na"
$ws.Range("G17").Value = "na"
$ws.Range("G18").Value = "na"
$ws.Range("G19").Value = "na"
$ws.Range("G20").Value = "na"
$ws.Range("G21").Value = "na"
$ws.Range("G22").Value = "na"
$ws.Range("G23").Value = "na"
$ws.Range("G24").Value = "na"
$ws.Range("G25").Value = "na"
$ws.Range("G26").Value = "na"
$ws.Range("G27").Value = "na"
$ws.Range("G28").Value = "na"
$ws.Range("G29").Value = "na"
$ws.Range("G30").Value = "na"
$ws.Range("G31").Value = "na"
$ws.Range("G32").Value = "na"
$ws.Range("G33").Value = "na"
$ws.Range("G34").Value = "na"
$ws.Range("G35").Value = "na"
$ws.Range("G36").Value = "na"
$ws.Range("G37").Value = "na"
$ws.Range("G38").Value = "na"
$ws.Range("G39").Value = "na"
$ws.Range("G40").Value = "na"
$ws.Range("G41").Value = "na"
$ws.Range("G42").Value = "na"
$ws.Range("G43").Value = "na"
$ws.Range("G44").Value = "na"
$ws.Range("G45").Value = "na"
$ws.Range("G46").Value = "na"
$ws.Range("G47").Value = "na"
$ws.Range("G48").Value = "na"
$ws.Range("G49").Value = "na"
$ws.Range("G50").Value = "na"
$ws.Range("G51").Value = "na"
$ws.Range("G52").Value = "na"
$ws.Range("G53").Value = "na"
$ws.Range("G54").Value = "na"
$ws.Range("G55").Value = "na"
$ws.Range("G56").Value = "na"
$ws.Range("G57").Value = "na"
$ws.Range("G58").Value = "na"
$ws.Range("G59").Value = "na"
$ws.Range("G60").Value = "na"
$ws.Range("G61").Value = "na"
$ws.Range("G62").Value = "na"
$ws.Range("G63").Value = "na"
$ws.Range("G64").Value = "na"
$ws.Range("G65").Value = "na"
$ws.Range("G66").Value = "na"
$ws.Range("G67").Value = "na"
$ws.Range("G68").Value = "na"
$ws.Range("G69").Value = "na"
$ws.Range("G70").Value = "na"
$ws.Range("G71").Value = "na"
$ws.Range("G72").Value = "na"
$ws.Range("G73").Value = "na"
$ws.Range("G74").Value = "na"
$ws.Range("G75").Value = "na"
$ws.Range("G76").Value = "na"
$ws.Range("G77").Value = "na"
$ws.Range("G78").Value = "na"
$ws.Range("G79").Value = "na"
$ws.Range("G80").Value = "na"
$ws.Range("G81").Value = "na"
$ws.Range("G82").Value = "na"
$ws.Range("G83").Value = "na"
$ws.Range("G84").Value = "na"
$ws.Range("G85").Value = "na"
$ws.Range("G86").Value = "na"
$ws.Range("G87").Value = "na"
$ws.Range("G88").Value = "運べば"

# Leave the selection on the last cell that was filled in, like the author did
$ws.Range("G88").Select()
